# Recreate baseline and match dac scenario pars
#
# - CAP / CAP_NEW / INVESTMENT sheets: the year header row shifts from
#   2015..2110 (15 cols, B:P) to 2025..2110 (13 cols, B:N) - the first two
#   periods (2015, 2020) are dropped and columns O:P are removed entirely.
# - All four sheets (CAP / CAP_NEW / INVESTMENT / REMOVAL) have their whole
#   data body (every numeric cell outside the node_loc label column and the
#   header row) reset to 0 - the new baseline run produced all-zero output
#   for this report.

$wb = $excel.ActiveWorkbook

$newYears = @(2025, 2030, 2035, 2040, 2045, 2050, 2055, 2060, 2070, 2080, 2090, 2100, 2110)

# Sheets whose header/columns need to be recreated (2015-2110 -> 2025-2110).
$resizedSheets = @("CAP", "CAP_NEW", "INVESTMENT")

foreach ($sheetName in $resizedSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Zero out the whole data body first (still 15 columns wide at this point).
    $ws.Range("B2:P13").Value = 0

    # Rewrite the year header (B1:N1) with the new period set.
    for ($i = 0; $i -lt $newYears.Length; $i++) {
        $ws.Cells.Item(1, 2 + $i).Value = $newYears[$i]
    }

    # Drop the now-unused trailing columns (old 2100/2110 slots, O:P).
    $ws.Range("O1:P13").Delete()
}

# REMOVAL already has the 2025-2110 header/shape; only the values reset.
$wsRemoval = $wb.Worksheets.Item("REMOVAL")
$wsRemoval.Range("B2:N13").Value = 0
